$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new ticker row at the bottom of the list (A76)
$ws.Range("A76").Value = "GRT-USD"
